# Generate Report for Handoff
# The two tracked files (3314774a-... and 81f122b5-...) swap places in every
# status table, and 3314774a-... picks up a fresh "Ready for handoff" status
# (new handoff/handback timestamps + a stale-handback error) while
# 81f122b5-... keeps its previous "Handed back: in sync with en-US" data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "81f122b5-5586-4a39-8579-09d726ed7697.md"
$ov.Range("B2").Value = "e2e\81f122b5-5586-4a39-8579-09d726ed7697.md"

$ov.Range("A3").Value = "3314774a-5bcf-4947-8455-15fab0d08c42.md"
$ov.Range("B3").Value = "e2e\3314774a-5bcf-4947-8455-15fab0d08c42.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-24 16:49:42"

foreach ($hl in $ov.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\81f122b5-5586-4a39-8579-09d726ed7697.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\3314774a-5bcf-4947-8455-15fab0d08c42.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "81f122b5-5586-4a39-8579-09d726ed7697.md"
$zh.Range("G2").Value = "81f122b5-5586-4a39-8579-09d726ed7697.9905165af3c639b33cd0f76180e6d083c20e75e7.zh-cn.xlf"
$zh.Range("I2").Value = "81f122b5-5586-4a39-8579-09d726ed7697.md"
$zh.Range("J2").Value = "81f122b5-5586-4a39-8579-09d726ed7697.9905165af3c639b33cd0f76180e6d083c20e75e7.zh-cn.xlf"

$zh.Range("A3").Value = "3314774a-5bcf-4947-8455-15fab0d08c42.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "3314774a-5bcf-4947-8455-15fab0d08c42.a1084e96799d82e3af0fdf3749bef954eea0f964.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-24 16:49:36"
$zh.Range("I3").Value = "3314774a-5bcf-4947-8455-15fab0d08c42.md"
$zh.Range("J3").Value = "3314774a-5bcf-4947-8455-15fab0d08c42.a1084e96799d82e3af0fdf3749bef954eea0f964.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/43d4c58795f5f35ee937e853d65d200e649a7fc9/e2e/3314774a-5bcf-4947-8455-15fab0d08c42.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6df4164a7f4a80d013912003fdd6defdb7f361d8/e2e/3314774a-5bcf-4947-8455-15fab0d08c42.md."

foreach ($hl in $zh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "81f122b5-5586-4a39-8579-09d726ed7697.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "81f122b5-5586-4a39-8579-09d726ed7697.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "3314774a-5bcf-4947-8455-15fab0d08c42.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "3314774a-5bcf-4947-8455-15fab0d08c42.md"
    }
}

$zh.Columns.Item(16).ColumnWidth = 39.1666667

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "81f122b5-5586-4a39-8579-09d726ed7697.md"
$de.Range("G2").Value = "81f122b5-5586-4a39-8579-09d726ed7697.9905165af3c639b33cd0f76180e6d083c20e75e7.de-de.xlf"
$de.Range("I2").Value = "81f122b5-5586-4a39-8579-09d726ed7697.md"
$de.Range("J2").Value = "81f122b5-5586-4a39-8579-09d726ed7697.9905165af3c639b33cd0f76180e6d083c20e75e7.de-de.xlf"

$de.Range("A3").Value = "3314774a-5bcf-4947-8455-15fab0d08c42.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "3314774a-5bcf-4947-8455-15fab0d08c42.a1084e96799d82e3af0fdf3749bef954eea0f964.de-de.xlf"
$de.Range("H3").Value = "2016-08-24 16:49:42"
$de.Range("I3").Value = "3314774a-5bcf-4947-8455-15fab0d08c42.md"
$de.Range("J3").Value = "3314774a-5bcf-4947-8455-15fab0d08c42.a1084e96799d82e3af0fdf3749bef954eea0f964.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/43d4c58795f5f35ee937e853d65d200e649a7fc9/e2e/3314774a-5bcf-4947-8455-15fab0d08c42.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6df4164a7f4a80d013912003fdd6defdb7f361d8/e2e/3314774a-5bcf-4947-8455-15fab0d08c42.md."

foreach ($hl in $de.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "81f122b5-5586-4a39-8579-09d726ed7697.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "81f122b5-5586-4a39-8579-09d726ed7697.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "3314774a-5bcf-4947-8455-15fab0d08c42.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "3314774a-5bcf-4947-8455-15fab0d08c42.md"
    }
}

$de.Columns.Item(16).ColumnWidth = 39.1666667
